$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update the Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-05T14:17:01+00:00"

# --- Sheet "Elements": update Binding Strength and Binding Value Set ---
$elements = $wb.Worksheets.Item("Elements")

# Row 5 / Row 6 - Binding Strength column X: required -> preferred
$elements.Range("X5").Value = "preferred"
$elements.Range("X6").Value = "preferred"

# Row 6 - Binding Value Set column Z: updated JDV URL
$elements.Range("Z6").Value = "https://mos.esante.gouv.fr/NOS/JDV_J01-XdsAuthorSpecialty-CISIS/FHIR/JDV-J01-XdsAuthorSpecialty-CISIS"

# Column Z width grows to fit the new (longer) URL text.
# (Target stored width is 83.80078125; the host's column-width model
# quantizes to whole pixels, so 83.0 is the nearest input that round-trips
# to the closest achievable stored width.)
$elements.Columns.Item(26).ColumnWidth = 83.0
